$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Epic summary data: Epic, PASSED, FAILED, BROKEN, SKIPPED, totalTests, passRate, status
$data = @(
    ,@('INVHUB-1018: User/System Preferences', 1, 0, 0, 0, 1, 100, 'Acceptable')
    ,@('INVHUB-1035: Sensa Copilot for Investigators GA', 2, 0, 0, 0, 2, 100, 'Acceptable')
    ,@('INVHUB-10821: Filtering on transactions', 8, 0, 0, 0, 8, 100, 'Acceptable')
    ,@('INVHUB-1180: Sensa Copilot enhancements/support(Q1)', 2, 0, 0, 0, 2, 100, 'Acceptable')
    ,@('INVHUB-13728: UI doc summary - Enterprise Knowledge Management', 4, 0, 0, 0, 4, 100, 'Acceptable')
    ,@('INVHUB-13905: (Phase 2) Workflow Actions mgmt', 28, 0, 0, 0, 28, 100, 'Acceptable')
    ,@('INVHUB-18: Authentication and Login', 5, 0, 0, 0, 5, 100, 'Acceptable')
    ,@('INVHUB-4521: NFR - Security', 2, 0, 0, 0, 2, 100, 'Acceptable')
    ,@('INVHUB-5339: CDS UI Component Library - Migration', 1, 0, 0, 0, 1, 100, 'Acceptable')
    ,@('INVHUB-7539: Tags Management', 1, 0, 0, 0, 1, 100, 'Acceptable')
    ,@('INVHUB-8565: Refinements to event/audit log', 1, 0, 0, 0, 1, 100, 'Acceptable')
    ,@('INVHUB-9173: Abanca Technical POC', 2, 0, 0, 0, 2, 100, 'Acceptable')
    ,@('INVHUB-9920: Tag Category (workflow)', 12, 0, 0, 0, 12, 100, 'Acceptable')
    ,@('NVHUB-5339: CDS UI Component Library - Migration', 3, 0, 0, 0, 3, 100, 'Acceptable')
    ,@('INVHUB-9369: IH : Language dropdown to generate Narrative  - No EPIC Tagged', 2, 0, 0, 0, 2, 100, 'Acceptable')
    ,@('INVHUB-2377: Roles Management Q1', 31, 1, 0, 0, 32, 96.88, 'Acceptable')
    ,@('INVHUB-9919: (Phase 3) Workflow - Decision making and multi level workflows', 52, 0, 2, 0, 54, 96.3, 'Acceptable')
    ,@('API Testing: GraphQL', 163, 6, 1, 0, 170, 95.88, 'Acceptable')
    ,@('INVHUB-3309: Investigation Hub Tags', 28, 2, 0, 0, 30, 93.33, 'Maintenance Advised')
    ,@('INVHUB-2452: Enhancements to the investigator experience', 193, 0, 0, 21, 214, 90.19, 'Maintenance Advised')
    ,@('INVHUB-2168: Teams Management - Q1', 7, 1, 0, 0, 8, 87.5, 'Maintenance Advised')
    ,@('INVHUB-3120: Workflows -Manager creation of Actions / map to teams', 81, 15, 0, 0, 96, 84.38, 'Maintenance Advised')
    ,@('INVHUB-9920: Maintenance of Tags and Tag Categories', 32, 6, 0, 0, 38, 84.20999999999999, 'Maintenance Advised')
    ,@('API Testing', 200, 24, 6, 9, 239, 83.68000000000001, 'Maintenance Advised')
    ,@('INVHUB-8981: Narrative contribution and UI view', 3, 1, 0, 0, 4, 75, 'Review Required')
    ,@('INVHUB-8114: Authentication', 30, 7, 4, 0, 41, 73.17, 'Review Required')
    ,@('INVHUB-126: Entity Centric Investigation (Subject Details) (MVP)', 30, 14, 1, 1, 46, 65.22, 'Review Required')
    ,@('INVHUB-1047: Subject Centric Investigation - TM', 5, 1, 2, 0, 8, 62.5, 'Review Required')
    ,@('INVHUB-11811: Prefill transactions on the transaction tab', 12, 11, 0, 0, 23, 52.17, 'Review Required')
    ,@('INVHUB-157: Assign (Reassign) and drill down (GA)', 1, 1, 0, 0, 2, 50, 'Review Required')
    ,@('INVHUB-2136: Subject Centric Investigation(Q1)', 1, 1, 0, 0, 2, 50, 'Review Required')
    ,@('INVHUB-11621: Provide detailed analysis at the individual detection level', 27, 7, 2, 21, 57, 47.37, 'Review Required')
    ,@('Test Cases Not Tagged', 74, 31, 46, 14, 165, 44.85, 'Review Required')
    ,@('INVHUB-290: Teams Management', 7, 9, 0, 0, 16, 43.75, 'Review Required')
    ,@('INVHUB-13: Menu and User Administration', 4, 3, 3, 0, 10, 40, 'Review Required')
    ,@('INVHUB-11963: Copilot: Add  some  more complex and standard  queries for Transactions', 5, 8, 0, 0, 13, 38.46, 'Review Required')
    ,@('INVHUB-7367: IH Bugs', 1, 2, 5, 0, 8, 12.5, 'Review Required')
    ,@('INVHUB-4976: Data Access - Q2', 4, 0, 32, 0, 36, 11.11, 'Review Required')
    ,@('INVHUB-13708: Transaction Ingestion Phase 1', 1, 9, 0, 0, 10, 10, 'Review Required')
    ,@('INVHUB-171: Data Access - GA', 0, 1, 15, 0, 16, 0, 'Review Required')
    ,@('INVHUB-5: Open Investigations List', 0, 4, 0, 0, 4, 0, 'Review Required')
    ,@('INVHUB-966: Open Investigations queue/list (GA)', 0, 4, 0, 0, 4, 0, 'Review Required')
    ,@('INVHUB-135: Dual browser capability - No EPIC Tagged', 0, 0, 3, 0, 3, 0, 'Review Required')
    ,@('INVHUB-8203: Remove Sensa Roles: Investigation Hub Roles to become the only Roles available - No EPIC Tagged', 0, 1, 0, 0, 1, 0, 'Review Required')
    ,@('TOTAL', 1066, 170, 122, 66, 1424, $null, $null)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    if ($null -eq $row[6]) {
        $ws.Cells.Item($rowIndex, 7).Value = ""
    } else {
        $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    }
    if ($null -eq $row[7]) {
        $ws.Cells.Item($rowIndex, 8).Value = ""
    } else {
        $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    }
    $rowIndex++
}

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
